$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-21: rewrite FirstName/LastName/Email/Phone/Year per the updated
# resource data. The Id column (A) stays a plain 1..N sequence and is
# untouched except for the newly-appended row.

# Note: E4's phone value is purely numeric ("123987654"); prefix it with an
# apostrophe so Excel keeps storing it as text instead of coercing it to a
# number (matches the source workbook's t="str" cell). E14's phone becomes
# blank; a lone apostrophe yields an empty *text* cell (t="str", empty <v/>)
# instead of clearing the cell to a true empty/number cell.
$rows = @(
    @{ r = 4;  b = "User6";   c = "User6LN";  d = "user6@gmail.com";   e = "'123987654";       f = 2000 }
    @{ r = 5;  b = "User4";   c = "User4LN";  d = "user4@gmail.com";   e = "(123)123456789";   f = 2000 }
    @{ r = 6;  b = "Boxis";   c = "Strong";   d = "boxis@gmail.com";   e = "(111)333222";       f = 1993 }
    @{ r = 7;  b = "Tor";     c = "Asgaard";  d = "tor@gmail.com";     e = "(111)333888";       f = 1994 }
    @{ r = 8;  b = "User1";   c = "Admin1";   d = "user1@gmail.com";   e = "(123)456780";       f = 1990 }
    @{ r = 9;  b = "Gunnar";  c = "Jensen";   d = "gunnar@gmail.com";  e = "(111)222444";       f = 1980 }
    @{ r = 10; b = "Bruce";   c = "Lee";      d = "bruce@gmail.com";   e = "(111)333445";       f = 1987 }
    @{ r = 11; b = "Gamora";  c = "Gamorak";  d = "gamora@gmail.com";  e = "(111)333111";       f = 1988 }
    @{ r = 12; b = "Witcher"; c = "Moon";     d = "witcher@gmail.com"; e = "(111)333999";       f = 1990 }
    @{ r = 13; b = "Supwom";  c = "Nanual";   d = "supwom@gmail.com";  e = "(111)333777";       f = 1988 }
    @{ r = 14; b = "User5";   c = "User5LN";  d = "user5@gmail.com";   e = "'";                 f = 0 }
    @{ r = 16; b = "User2";   c = "User2LN";  d = "user2@gmail.com";   e = "(123)123123";       f = 1995 }
    @{ r = 17; b = "Lee";     c = "Christmas";d = "lee@gmail.com";     e = "(111)333444";       f = 1977 }
    @{ r = 18; b = "User3";   c = "User3LN";  d = "user3@gmail.com";   e = "(123)123124";       f = 1996 }
    @{ r = 19; b = "Marvel";  c = "Levram";   d = "marvel@gmail.com";  e = "(111)333555";       f = 1995 }
    @{ r = 20; b = "Jean";    c = "Vilain";   d = "jean@gmail.com";    e = "(111)222777";       f = 1973 }
    @{ r = 21; b = "Sonya";   c = "Night";    d = "sonya@gmail.com";   e = "(111)333666";       f = 1996 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
    $ws.Cells.Item($r, 6).Value = $row.f
}

# Row 22 is brand new (Natalia Romanoff, formerly at row 21, now pushed down
# a row). Clone row 21's formatting onto row 22 first so the new row matches
# the rest of the table, then overwrite it with Natalia's data.
$ws.Range("A21:F21").Copy($ws.Range("A22:F22"))
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "Natalia"
$ws.Cells.Item(22, 3).Value = "Romanoff"
$ws.Cells.Item(22, 4).Value = "natalia@gmail.com"
$ws.Cells.Item(22, 5).Value = "(111)222888"
$ws.Cells.Item(22, 6).Value = 1986
